$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHD")

# Row 4 - Inventory
$ws.Range("B4").Value = 541000000.0
$ws.Range("C4").Value = 495000000.0
$ws.Range("D4").Value = 498000000.0
$ws.Range("E4").Value = 456000000.0
$ws.Range("F4").Value = 396000000.0

# Row 13 - Accounts Payable
$ws.Range("C13").Value = 588000000.0
$ws.Range("D13").Value = 581000000.0
$ws.Range("E13").Value = 517000000.0
$ws.Range("F13").Value = 460000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("C20").Value = 707000000.0
$ws.Range("D20").Value = 600000000.0
$ws.Range("E20").Value = 582000000.0
$ws.Range("F20").Value = 580000000.0

# Row 33 - Net Debt
$ws.Range("G33").Value = 1907400000.0

# Row 34 - Total Debt
$ws.Range("G34").Value = 2063100000.0
